$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.738.79'
$ws.Range('E2').Value = '  -1.19%  '

$ws.Range('D3').Value = '2.102.25'
$ws.Range('E3').Value = '  -0.17%  '

$ws.Range('E4').Value = '  +0.57%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '346.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.68%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5203'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.68%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4403'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.21%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.93'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.63%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09404'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.172'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.42%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.06%  '

$ws.Range('D13').Value = '2.100.46'
$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.822'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.227'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.77%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '102.84'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.39%  '

$ws.Range('E17').Value = '  +0.62%  '


$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.36%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06671'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '

$ws.Range('E21').Value = '  +0.49%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.267'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.25%  '

$ws.Range('D23').Value = '29.796.77'
$ws.Range('E23').Value = '  -1.32%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.44%  '

$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.427.22'
$ws.Range('E25').Value = '  +3.10%  '

$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.326'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.00%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.515'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.137'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.76%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.720'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.93%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1054'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.18%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.218'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.55%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.951'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.15%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.350'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.81%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.50'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.08%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02586'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.17%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06743'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.18%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7005'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.56%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.56%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.335'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.01%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2222'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.70%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6841'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.84%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.47'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.39%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.356'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.27%  '

$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000358'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.82%  '

$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.635'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.41%  '

$ws.Range('B49').Value = 'WEMIXTOKEN'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.218'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.51%  '

$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.220'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.16%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.52%  '
